# ---------------------------------------------------------------------------
# Fix bug in comparison of GT tables
#
# 1) Note the bugfix on row 78 (test 5 GS block)
# 2) Fill in the "RW" (O:Y) sub-block for test 5 (rows 79-87), mirroring the
#    pattern already present for test 4 (rows 69-77)
# 3) Note re-running after the fix on row 97
# 4) Add a new "test 7" GS block (rows 98-109), re-comparing after the fix
# 5) Leave the selection on the newly added data
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bugfix note on row 78 -------------------------------------------------
$ws.Range("P78").Value = "fix features returning NA for st.error measurements"

# 2) Test 5 "RW" sub-block (O:Y), rows 79-87 -------------------------------
$rwBlock = @(
  @(79, 5, "RW", "rf BS15_AU_02a_files_1-104", 581, 268, 313, 12, 0.957142857142857, 0.538726333907057, 0.856230031948882, 0.96441925118333),
  @(80, 5, "RW", "rf BS14_AU_04_files_All", 1408, 410, 998, 51, 0.889370932754881, 0.708806818181818, 0.410821643286573, 0.96441925118333),
  @(81, 5, "RW", "rf AW12_AU_BS3_files_All", 1907, 771, 1136, 82, 0.903868698710434, 0.595700052438385, 0.678697183098592, 0.96441925118333),
  @(82, 5, "RW", "rf BS13_AU_04_files_All", 3278, 1376, 1902, 228, 0.85785536159601, 0.580231848688225, 0.723449001051525, 0.96441925118333),
  @(83, 5, "RW", "rf BS16_AU_02a_files_1-175", 1089, 453, 636, 37, 0.924489795918367, 0.584022038567493, 0.712264150943396, 0.96441925118333),
  @(84, 5, "RW", "rf BS15_AU_02b_files_All", 508, 144, 364, 33, 0.813559322033898, 0.716535433070866, 0.395604395604396, 0.96441925118333),
  @(85, 5, "RW", "rf AW14_AU_BS3_files_1-160", 1569, 541, 1028, 65, 0.892739273927393, 0.655194391332059, 0.526264591439689, 0.96441925118333),
  @(86, 5, "RW", "rf AL16_AU_BS1_files_All", 1456, 409, 1047, 30, 0.931662870159453, 0.719093406593407, 0.390639923591213, 0.96441925118333),
  @(87, 5, "RW", "rf all", 11796, 4372, 7424, 538, 0.890427698574338, 0.629365886741268, 0.588900862068966, 0.96441925118333)
)

foreach ($row in $rwBlock) {
  $r = $row[0]
  for ($i = 1; $i -lt $row.Count; $i++) {
    $col = 14 + $i   # data starts at column O (15)
    $ws.Cells.Item($r, $col).Value = $row[$i]
  }
}

# 3) Note that data below is after the fix ---------------------------------
$ws.Range("C97").Value = "after fix adaptive compare"

# 4) Test 7 "GS" block (A:K), rows 98-109 ----------------------------------
$gsBlock = @(
  @(98, 7, "GS", "rf AW15_AU_BS3_files_705-749", 1014, 653, 361, 44, 0.93687230989957, 0.356015779092702, 1.80886426592798, 0.948700021365922),
  @(99, 7, "GS", "rf BS12_AU_02a_files_1-46", 537, 401, 136, 0, 1, 0.253258845437616, 2.94852941176471, 0.948700021365922),
  @(100, 7, "GS", "rf AW14_AU_BS3_files_1-71", 1498, 965, 533, 54, 0.947006869479882, 0.355807743658211, 1.81050656660413, 0.948700021365922),
  @(101, 7, "GS", "rf BS13_AU_04_files_137-224", 1460, 685, 775, 4, 0.994194484760523, 0.530821917808219, 0.883870967741936, 0.948700021365922),
  @(102, 7, "GS", "rf AW12_AU_BS3_files_1-250", 2987, 1695, 1292, 112, 0.938018815716657, 0.432541011047874, 1.31191950464396, 0.948700021365922),
  @(103, 7, "GS", "rf AW12_AU_BS3_files_1464-1507", 778, 562, 216, 57, 0.907915993537964, 0.277634961439589, 2.60185185185185, 0.948700021365922),
  @(104, 7, "GS", "rf AW14_AU_BS3_files_309-369", 865, 622, 243, 96, 0.866295264623955, 0.28092485549133, 2.559670781893, 0.948700021365922),
  @(105, 7, "GS", "rf AW15_AU_BS2_files_33-103", 745, 407, 338, 67, 0.858649789029536, 0.453691275167785, 1.20414201183432, 0.948700021365922),
  @(106, 7, "GS", "rf AL16_AU_BS3_files_77-170", 629, 319, 310, 47, 0.871584699453552, 0.492845786963434, 1.02903225806452, 0.948700021365922),
  @(107, 7, "GS", "rf BS12_AU_02b_files_689-747", 1986, 1472, 514, 241, 0.859311150029189, 0.258811681772407, 2.86381322957198, 0.948700021365922),
  @(108, 7, "GS", "rf BS14_AU_04_files_74-148", 674, 343, 331, 57, 0.8575, 0.491097922848665, 1.03625377643505, 0.948700021365922),
  @(109, 7, "GS", "rf all", 13173, 8124, 5049, 779, 0.912501404021116, 0.383283989979504, 1.60903149138443, 0.948700021365922)
)

foreach ($row in $gsBlock) {
  $r = $row[0]
  for ($i = 1; $i -lt $row.Count; $i++) {
    $col = $i        # data starts at column A (1)
    $ws.Cells.Item($r, $col).Value = $row[$i]
  }
}

# 5) Move the selection onto the newly entered data ------------------------
$ws.Range("C98").Select()
